$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 367.5
$ws.Range("I2").Value = 118.8
$ws.Range("J2").Value = 616.2
$ws.Range("K2").Value = 118.8
$ws.Range("L2").Value = 616.2
$ws.Range("M2").Value = -5.799999999999997
$ws.Range("N2").Value = -842.2
$ws.Range("H9").Value = 81.21429000000001
$ws.Range("I9").Value = 57.333332
$ws.Range("J9").Value = 124.2
$ws.Range("K9").Value = 57.333332
$ws.Range("L9").Value = 124.2
$ws.Range("M9").Value = 111.666668
$ws.Range("N9").Value = -462.2
$ws.Range("H28").Value = 672.9167
$ws.Range("I28").Value = 638.5714
$ws.Range("J28").Value = 721
$ws.Range("K28").Value = 638.5714
$ws.Range("L28").Value = 721
$ws.Range("M28").Value = -153.5714
$ws.Range("N28").Value = -1691
$ws.Range("H38").Value = 594.2353000000001
$ws.Range("I38").Value = 140.875
$ws.Range("J38").Value = 997.2222
$ws.Range("K38").Value = 422.625
$ws.Range("L38").Value = 2991.6666
$ws.Range("M38").Value = -50.625
$ws.Range("N38").Value = -3735.6666
$ws.Range("H40").Value = 2305.8125
$ws.Range("I40").Value = 1987.5
$ws.Range("J40").Value = 2624.125
$ws.Range("K40").Value = 1987.5
$ws.Range("L40").Value = 2624.125
$ws.Range("M40").Value = -1812.5
$ws.Range("N40").Value = -2974.125
$ws.Range("H58").Value = 1694.6522
$ws.Range("I58").Value = 343.07693
$ws.Range("J58").Value = 3451.7
$ws.Range("K58").Value = 1029.23079
$ws.Range("L58").Value = 10355.1
$ws.Range("M58").Value = -879.2307900000001
$ws.Range("N58").Value = -10655.1
$ws.Range("H74").Value = 6228.8335
$ws.Range("I74").Value = 6167.8184
$ws.Range("J74").Value = 6900
$ws.Range("K74").Value = 6167.8184
$ws.Range("L74").Value = 6900
$ws.Range("M74").Value = -5231.8184
$ws.Range("N74").Value = -8772
$ws.Range("H76").Value = 4606.5557
$ws.Range("I76").Value = 3994.8333
$ws.Range("K76").Value = 3994.8333
$ws.Range("M76").Value = -3679.8333
$ws.Range("H77").Value = 6228.8335
$ws.Range("I77").Value = 6167.8184
$ws.Range("J77").Value = 6900
$ws.Range("K77").Value = 30839.092
$ws.Range("L77").Value = 34500
$ws.Range("M77").Value = -26159.092
$ws.Range("N77").Value = -43860
$ws.Range("H79").Value = 4606.5557
$ws.Range("I79").Value = 3994.8333
$ws.Range("K79").Value = 3994.8333
$ws.Range("M79").Value = -2902.8333
$ws.Range("H106").Value = 28999.375
$ws.Range("I106").Value = 25932.666
$ws.Range("K106").Value = 25932.666
$ws.Range("M106").Value = -25301.666
$ws.Range("H116").Value = 4033.375
$ws.Range("I116").Value = 2857.6
$ws.Range("J116").Value = 5993
$ws.Range("K116").Value = 2857.6
$ws.Range("L116").Value = 5993
$ws.Range("M116").Value = 584.4000000000001
$ws.Range("N116").Value = -12877
$ws.Range("H135").Value = 1193.5454
$ws.Range("J135").Value = 3333
$ws.Range("L135").Value = 29997
$ws.Range("N135").Value = -35067

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 3863.2
$ws.Range("I25").Value = 2105.3333
$ws.Range("K25").Value = 2105.3333
$ws.Range("M25").Value = -1703.3333
$ws.Range("H35").Value = 2436.25
$ws.Range("I35").Value = 2436.25
$ws.Range("K35").Value = 2436.25
$ws.Range("M35").Value = -2030.25
$ws.Range("H88").Value = 1449.8334
$ws.Range("I88").Value = 692
$ws.Range("J88").Value = 2207.6667
$ws.Range("K88").Value = 692
$ws.Range("L88").Value = 2207.6667
$ws.Range("M88").Value = -286
$ws.Range("N88").Value = -3019.6667
$ws.Range("H91").Value = 1449.8334
$ws.Range("I91").Value = 692
$ws.Range("J91").Value = 2207.6667
$ws.Range("K91").Value = 692
$ws.Range("L91").Value = 2207.6667
$ws.Range("M91").Value = 712
$ws.Range("N91").Value = -5015.6667
$ws.Range("H102").Value = 6948036
$ws.Range("I102").Value = 7411105
$ws.Range("K102").Value = 7411105
$ws.Range("M102").Value = -7409483

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3305.5715
$ws.Range("I20").Value = 4187.8
$ws.Range("J20").Value = 1100
$ws.Range("K20").Value = 4187.8
$ws.Range("L20").Value = 1100
$ws.Range("M20").Value = -3940.8
$ws.Range("N20").Value = -1594
$ws.Range("H86").Value = 1249.5
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 1249.5
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H94").Value = 1683.75
$ws.Range("I94").Value = 828.4167
$ws.Range("K94").Value = 828.4167
$ws.Range("M94").Value = -377.4167
$ws.Range("H99").Value = 863
$ws.Range("I99").Value = 720.875
$ws.Range("K99").Value = 720.875
$ws.Range("M99").Value = 777.125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 47999
$ws.Range("H62").Value = 133766.33
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 200149.5
$ws.Range("K62").Value = 1000
$ws.Range("L62").Value = 200149.5
$ws.Range("M62").Value = -376
$ws.Range("N62").Value = -201397.5
$ws.Range("H65").Value = 133766.33
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 200149.5
$ws.Range("K65").Value = 5000
$ws.Range("L65").Value = 1000747.5
$ws.Range("M65").Value = -1880
$ws.Range("N65").Value = -1006987.5
$ws.Range("H99").Value = 12042.549
$ws.Range("I99").Value = 6314.8237
$ws.Range("K99").Value = 6314.8237
$ws.Range("M99").Value = -4816.8237
$ws.Range("H126").Value = 12042.549
$ws.Range("I126").Value = 6314.8237
$ws.Range("K126").Value = 18944.4711
$ws.Range("M126").Value = -16474.4711
$ws.Range("H134").Value = 3571.6667
$ws.Range("I134").Value = 3473.5
$ws.Range("J134").Value = 3620.75
$ws.Range("K134").Value = 10420.5
$ws.Range("L134").Value = 10862.25
$ws.Range("M134").Value = -7885.5
$ws.Range("N134").Value = -15932.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1426.6666
$ws.Range("I32").Value = 1666.6666
$ws.Range("J32").Value = 946.6667
$ws.Range("K32").Value = 4999.9998
$ws.Range("L32").Value = 2840.0001
$ws.Range("M32").Value = -4716.9998
$ws.Range("N32").Value = -3406.0001
$ws.Range("H120").Value = 11333
$ws.Range("I120").Value = 1999
$ws.Range("K120").Value = 5997
$ws.Range("M120").Value = -1159

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15152.385
$ws.Range("I70").Value = 11622.375
$ws.Range("K70").Value = 11622.375
$ws.Range("M70").Value = -11352.375
$ws.Range("H73").Value = 15152.385
$ws.Range("I73").Value = 11622.375
$ws.Range("K73").Value = 11622.375
$ws.Range("M73").Value = -10686.375
$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492
$ws.Range("H126").Value = 5699.8184
$ws.Range("I126").Value = 5549.857
$ws.Range("J126").Value = 5962.25
$ws.Range("K126").Value = 16649.571
$ws.Range("L126").Value = 17886.75
$ws.Range("M126").Value = -14179.571
$ws.Range("N126").Value = -22826.75
$ws.Range("H129").Value = 45000
$ws.Range("J129").Value = 45000
$ws.Range("L129").Value = 45000
$ws.Range("N129").Value = -55000

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 4799.8
$ws.Range("J3").Value = 4799.8
$ws.Range("L3").Value = 4799.8
$ws.Range("N3").Value = -5023.8
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H15").Value = 4799.8
$ws.Range("J15").Value = 4799.8
$ws.Range("L15").Value = 4799.8
$ws.Range("N15").Value = -5139.8
$ws.Range("H46").Value = 4011.0454
$ws.Range("J46").Value = 6332.5557
$ws.Range("L46").Value = 6332.5557
$ws.Range("N46").Value = -6708.5557
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H68").Value = 2850
$ws.Range("I68").Value = 2850
$ws.Range("K68").Value = 2850
$ws.Range("M68").Value = -2101
$ws.Range("H71").Value = 2850
$ws.Range("I71").Value = 2850
$ws.Range("K71").Value = 14250
$ws.Range("M71").Value = -10506
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("N106").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 24000
$ws.Range("I9").Value = 24000
$ws.Range("K9").Value = 24000
$ws.Range("M9").Value = -23860
$ws.Range("H14").Value = 25000
$ws.Range("J14").Value = 25000
$ws.Range("L14").Value = 25000
$ws.Range("N14").Value = -25336
$ws.Range("H40").Value = 34312.375
$ws.Range("I40").Value = 32000
$ws.Range("J40").Value = 34642.715
$ws.Range("K40").Value = 32000
$ws.Range("L40").Value = 34642.715
$ws.Range("M40").Value = -31851
$ws.Range("N40").Value = -34940.715
$ws.Range("H47").Value = 21332.666
$ws.Range("J47").Value = 21332.666
$ws.Range("L47").Value = 21332.666
$ws.Range("N47").Value = -22476.666
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H136").Value = 3448
$ws.Range("I136").Value = 3439.4167
$ws.Range("K136").Value = 10318.2501
$ws.Range("M136").Value = -7768.250100000001

Write-Host "Applied all changes"